$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: C3,D3,E3,F3
$ws.Range("C3").Value = "39"
$ws.Range("D3").Value = "27"
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "2"

# Row 5: C5,F5 (D5,E5 unchanged)
$ws.Range("C5").Value = "30"
$ws.Range("F5").Value = "0"

# Row 6: C6,D6,E6 (F6 unchanged)
$ws.Range("C6").Value = "8"
$ws.Range("D6").Value = "9"
$ws.Range("E6").Value = "1"

# Row 7: C7,D7,E7,F7
$ws.Range("C7").Value = "41"
$ws.Range("D7").Value = "34"
$ws.Range("E7").Value = "3"
$ws.Range("F7").Value = "2"

# Row 8: C8,D8,E8,F8
$ws.Range("C8").Value = "45"
$ws.Range("D8").Value = "25"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "4"

# Row 9: C9,D9,E9,F9
$ws.Range("C9").Value = "38"
$ws.Range("D9").Value = "20"
$ws.Range("E9").Value = "5"
$ws.Range("F9").Value = "1"

# Row 10: C10,D10,E10,F10
$ws.Range("C10").Value = "71"
$ws.Range("D10").Value = "48"
$ws.Range("E10").Value = "6"
$ws.Range("F10").Value = "3"

# Row 11: C11,D11,E11,F11
$ws.Range("C11").Value = "42"
$ws.Range("D11").Value = "40"
$ws.Range("E11").Value = "4"
$ws.Range("F11").Value = "0"

# Row 12: C12,D12,E12 (F12 unchanged)
$ws.Range("C12").Value = "13"
$ws.Range("D12").Value = "19"
$ws.Range("E12").Value = "2"
